# SLR_raw_concerns.xlsx — reorganize "FirstOther" (shared-string 13) entries
# into their real concern categories, and refresh a couple of view-state bits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Re-label cells that were previously the generic "FirstOther" bucket ----
$ws.Range("B128").Value = "Haz Waste/Industry"
$ws.Range("B164").Value = "Fiscal Cost"
$ws.Range("B191").Value = "Other"
$ws.Range("B201").Value = "Housing"
$ws.Range("B240").Value = "Other"
$ws.Range("B263").Value = "Haz Waste/Industry"
$ws.Range("B278").Value = "Fiscal Cost"
$ws.Range("B280").Value = "Other Infra"

$ws.Range("B313").Value = "DACs"
$ws.Range("C313").ClearContents() | Out-Null

$ws.Range("B377").Value = "Econ Growth"
$ws.Range("B396").Value = "Other Infra"
$ws.Range("B399").Value = "Stormwater"

$ws.Range("B418").Value = "Transpo"
$ws.Range("C418").ClearContents() | Out-Null

$ws.Range("B427").Value = "Flooding"
$ws.Range("B429").Value = "Other Infra"
$ws.Range("B459").Value = "Ecosystem"
$ws.Range("B473").Value = "Housing"
$ws.Range("C481").Value = "Commercial"
$ws.Range("C499").Value = "Haz Waste/Industry"
$ws.Range("C507").Value = "Cultural Resources"
$ws.Range("A512").Value = "Ecosystem"
$ws.Range("B554").Value = "DACs"

# Row 599 previously only had A/B populated; it now gains a C value too.
$ws.Range("C599").Value = "Haz Waste/Industry"

$ws.Range("B615").Value = "Water"
$ws.Range("C615").ClearContents() | Out-Null

$ws.Range("C666").Value = "Other Infra"
$ws.Range("A674").Value = "Flooding"
$ws.Range("B699").Value = "Other Infra"
$ws.Range("B701").Value = "Other"
$ws.Range("B703").Value = "Commercial"

# Row 730 previously had a stray B730 "FirstOther" cell with no A/C siblings worth keeping it for.
$ws.Range("B730").ClearContents() | Out-Null

$ws.Range("C748").Value = "Cultural Resources"
$ws.Range("B794").Value = "Other"
$ws.Range("B832").Value = "Ecosystem"

# ---- View-state refresh: selection moved from C13 to I6, top-left scroll reset ----
$ws.Range("I6").Select() | Out-Null

# ---- Workbook metadata: absPath now reflects the file's own subfolder ----
$wb.Path = "C:\Users\kyras\OneDrive\Desktop\SLRSurvey\SLRSurvey\CurrentData\Raw Edgelists_Excel"
